# Fix Order Filing Petition Template
#
# The document was re-saved by Word after a spell/grammar-check pass.
# For most of the edited paragraphs the *visible* text is unchanged but
# Word split existing runs and inserted <w:proofErr> bookmarks around the
# words/phrases it flagged (the "{{ ... }}" Jinja tags read as spelling /
# grammar issues). One paragraph also gained a genuine content fix: a
# missing space between "}}" and "to".
#
# We rebuild each affected paragraph's full Open XML (including the
# <w:proofErr/> markers and the new run boundaries) and splice it back in
# with Range.InsertXML, which replaces the whole paragraph range in place.

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Paragraph 1: "IN THE {{ court }} COURT OF {{ upper_county }} COUNTY, ..."
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("IN THE") | Out-Null
$p1 = $r1.Paragraphs(1).Range
$xml1 = '<w:p ' + $W + '>' +
    '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
    '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
    '<w:r><w:instrText xml:space="preserve"> SEQ CHAPTER \h \r 1</w:instrText></w:r>' +
    '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
    '<w:r><w:t>I</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">N THE </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>{{ court</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> COURT OF {{ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>upper_county</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> }} COUNTY, WEST VIRGINIA</w:t></w:r>' +
    '</w:p>'
$p1.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------
# Paragraph 6 (Title line): "{{ p.name.full(middle=...full...) }}"
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("p.name.full") | Out-Null
$p2 = $r2.Paragraphs(1).Range
$xml2 = '<w:p ' + $W + '>' +
    '<w:pPr><w:pStyle w:val="Title"/><w:jc w:val="left"/></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>p</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>.name.full</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(middle=' + [char]0x201D + 'full' + [char]0x201D + ') }}</w:t></w:r>' +
    '</w:p>'
$p2.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------
# Paragraph 11: "This day came ... Petitioner's name from {{ p.name.full
# (middle=...) }} to {{ new_name.name.full(middle=...) }}."
# (also fixes the missing space before "to")
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("This day came") | Out-Null
$p3 = $r3.Paragraphs(1).Range
$xml3 = '<w:p ' + $W + '>' +
    '<w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r>' +
    '<w:r><w:t xml:space="preserve">This day came the above-named Petitioner and tendered to the Court a petition, duly verified, requesting that this Court enter an ORDER changing the Petitioner' + [char]0x2019 + 's name from </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>p</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>.name.full</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">(middle=' + [char]0x201D + 'full' + [char]0x201D + ') </w:t></w:r>' +
    '<w:r><w:t>}}</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">to </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>new_name.name.full</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(middle=' + [char]0x201D + 'full' + [char]0x201D + ')</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>}}</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
$p3.InsertXML($xml3) | Out-Null

# ---------------------------------------------------------------------
# Paragraph 15: "ENTERED this ______ day of  _____________________, 20___."
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("ENTERED this") | Out-Null
$p4 = $r4.Paragraphs(1).Range
$xml4 = '<w:p ' + $W + '>' +
    '<w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">ENTERED this ______ day </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>of  _</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>____________________, 20___.</w:t></w:r>' +
    '</w:p>'
$p4.InsertXML($xml4) | Out-Null

Write-Output "done"
